$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New table data (rows 2-17): Player, Position, Team
$data = @(
    @("Tyler Herro",       "PG,SG",    "Miami Heat"),
    @("Luguentz Dort",     "SG,SF",    "Oklahoma City Thunder"),
    @("Evan Mobley",       "PF,C",     "Cleveland Cavaliers"),
    @("Miles Bridges",     "SF,PF",    "Charlotte Hornets"),
    @("Brook Lopez",       "C",        "Milwaukee Bucks"),
    @("Luka Doncic",       "PG,SG",    "Dallas Mavericks"),
    @("Scottie Barnes",    "SG,SF,PF", "Toronto Raptors"),
    @("Nikola Vucevic",    "PF,C",     "Chicago Bulls"),
    @("Mikal Bridges",     "SG,SF,PF", "New York Knicks"),
    @("DeMar DeRozan",     "SF,PF",    "Sacramento Kings"),
    @("Gradey Dick",       "SG,SF",    "Toronto Raptors"),
    @("De'Aaron Fox",      "PG",       "Sacramento Kings"),
    @("Ja Morant",         "PG",       "Memphis Grizzlies"),
    @("Santi Aldama",      "PF,C",     "Memphis Grizzlies"),
    @("Jonathan Kuminga",  "SF,PF",    "Golden State Warriors"),
    @("Josh Giddey",       "PG,SG,SF", "Chicago Bulls")
)

$row = 2
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $row++
}
